$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J header for year 2020 (same style as I3)
$ws.Range("J3").Value = 2020
$ws.Range("J3").Style = $ws.Range("I3").Style

# Row 4: hazardous waste generation value for 2020
$ws.Range("J4").Value = 11545.7
$ws.Range("J4").Style = $ws.Range("I4").Style

# Row 5: population value for 2020 -- stored as text "1 754,6"
$ws.Range("J5").Value = "1 754,6"
$ws.Range("J5").Style = $ws.Range("A5").Style

# Row 6: per-person value for 2020
$ws.Range("J6").Value = 6636.8
$ws.Range("J6").Style = $ws.Range("I6").Style

# Remove rows 27 and 28 (last two blank rows)
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()

$ws.Range("G22").Select()
